$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'315.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'2.32%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'39.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'-0.48%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Formula = "'5.140"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Formula = "'0.08183"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'0.68%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'1.980"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'2.01%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'4.368"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'3.36%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Formula = "'8.310"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "'2.12%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Formula = "'0.9386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "'1.05%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Formula = "'-8.45%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.1969"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'2.60%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'0.09032"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'-1.09%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Formula = "'0.03544"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Formula = "'0.83%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'0.09757"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'-0.51%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'0.001410"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'1.43%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'0.006252"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'5.14%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'3.635"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'-7.71%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Formula = "'-1.66%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Formula = "'1.22%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Formula = "'-2.35%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'4.960"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'6.83%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Formula = "'1.57%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'0.04370"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'0.06%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'0.001242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'1.10%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'0.004767"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'9.10%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'0.0003894"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'199.41%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Formula = "'-7.56%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Formula = "'0.02235"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'9.46%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'0.05204"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'2.96%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'0.007750"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'5.09%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Formula = "'0.01037"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'5.42%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'0.1402"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'2.66%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Formula = "'0.002102"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "'-1.35%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'0.008873"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'-5.33%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'0.00006827"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'7.13%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'0.08%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'0.003002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'10.03%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.001692"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'30.07%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'0.08%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'0.08%"
$ws.Range("E51").Style = "Normal"
